$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (Trimestre) to remain text so date-like strings are not
# auto-converted to date serials by Excel's input parser.
$dateRange = $ws.Range("C2:C64")
$dateRange.NumberFormat = "@"

$ws.Range("C2").Value = "01/10/2016"
$ws.Range("D2").Value = 2874
$ws.Range("C3").Value = "01/01/2017"
$ws.Range("D3").Value = 2910
$ws.Range("C4").Value = "01/04/2017"
$ws.Range("D4").Value = 2882
$ws.Range("C5").Value = "01/07/2017"
$ws.Range("D5").Value = 2893
$ws.Range("C6").Value = "01/10/2017"
$ws.Range("D6").Value = 2918
$ws.Range("C7").Value = "01/01/2018"
$ws.Range("D7").Value = 2939
$ws.Range("C8").Value = "01/04/2018"
$ws.Range("D8").Value = 2945
$ws.Range("C9").Value = "01/07/2018"
$ws.Range("D9").Value = 2934
$ws.Range("C10").Value = "01/10/2018"
$ws.Range("D10").Value = 2958
$ws.Range("C11").Value = "01/01/2019"
$ws.Range("D11").Value = 2974
$ws.Range("C12").Value = "01/04/2019"
$ws.Range("D12").Value = 2938
$ws.Range("C13").Value = "01/07/2019"
$ws.Range("D13").Value = 2943
$ws.Range("C14").Value = "01/10/2019"
$ws.Range("D14").Value = 2969
$ws.Range("C15").Value = "01/01/2020"
$ws.Range("D15").Value = 3003
$ws.Range("C16").Value = "01/04/2022"
$ws.Range("D16").Value = 2785
$ws.Range("C17").Value = "01/07/2022"
$ws.Range("D17").Value = 2887
$ws.Range("C18").Value = "01/10/2022"
$ws.Range("D18").Value = 2940
$ws.Range("C19").Value = "01/01/2023"
$ws.Range("D19").Value = 2959
$ws.Range("C20").Value = "01/04/2023"
$ws.Range("D20").Value = 2958
$ws.Range("C21").Value = "01/07/2023"
$ws.Range("D21").Value = 3007
$ws.Range("C22").Value = "01/10/2023"
$ws.Range("D22").Value = 3032
$ws.Range("C23").Value = "01/10/2016"
$ws.Range("D23").Value = 1938
$ws.Range("C24").Value = "01/01/2017"
$ws.Range("D24").Value = 1999
$ws.Range("C25").Value = "01/04/2017"
$ws.Range("D25").Value = 1997
$ws.Range("C26").Value = "01/07/2017"
$ws.Range("D26").Value = 1967
$ws.Range("C27").Value = "01/10/2017"
$ws.Range("D27").Value = 2032
$ws.Range("C28").Value = "01/01/2018"
$ws.Range("D28").Value = 2040
$ws.Range("C29").Value = "01/04/2018"
$ws.Range("D29").Value = 2037
$ws.Range("C30").Value = "01/07/2018"
$ws.Range("D30").Value = 2038
$ws.Range("C31").Value = "01/10/2018"
$ws.Range("D31").Value = 2054
$ws.Range("C32").Value = "01/01/2019"
$ws.Range("D32").Value = 2061
$ws.Range("C33").Value = "01/04/2019"
$ws.Range("D33").Value = 2038
$ws.Range("C34").Value = "01/07/2019"
$ws.Range("D34").Value = 2018
$ws.Range("C35").Value = "01/10/2019"
$ws.Range("D35").Value = 2048
$ws.Range("C36").Value = "01/01/2020"
$ws.Range("D36").Value = 2072
$ws.Range("C37").Value = "01/04/2022"
$ws.Range("D37").Value = 1864
$ws.Range("C38").Value = "01/07/2022"
$ws.Range("D38").Value = 1945
$ws.Range("C39").Value = "01/10/2022"
$ws.Range("D39").Value = 1961
$ws.Range("C40").Value = "01/01/2023"
$ws.Range("D40").Value = 2021
$ws.Range("C41").Value = "01/04/2023"
$ws.Range("D41").Value = 2004
$ws.Range("C42").Value = "01/07/2023"
$ws.Range("D42").Value = 2015
$ws.Range("C43").Value = "01/10/2023"
$ws.Range("D43").Value = 2040
$ws.Range("C44").Value = "01/10/2016"
$ws.Range("D44").Value = 2303
$ws.Range("C45").Value = "01/01/2017"
$ws.Range("D45").Value = 2317
$ws.Range("C46").Value = "01/04/2017"
$ws.Range("D46").Value = 2216
$ws.Range("C47").Value = "01/07/2017"
$ws.Range("D47").Value = 2101
$ws.Range("C48").Value = "01/10/2017"
$ws.Range("D48").Value = 2044
$ws.Range("C49").Value = "01/01/2018"
$ws.Range("D49").Value = 2146
$ws.Range("C50").Value = "01/04/2018"
$ws.Range("D50").Value = 2047
$ws.Range("C51").Value = "01/07/2018"
$ws.Range("D51").Value = 2157
$ws.Range("C52").Value = "01/10/2018"
$ws.Range("D52").Value = 2116
$ws.Range("C53").Value = "01/01/2019"
$ws.Range("D53").Value = 2111
$ws.Range("C54").Value = "01/04/2019"
$ws.Range("D54").Value = 2042
$ws.Range("C55").Value = "01/07/2019"
$ws.Range("D55").Value = 2013
$ws.Range("C56").Value = "01/10/2019"
$ws.Range("D56").Value = 1999
$ws.Range("C57").Value = "01/01/2020"
$ws.Range("D57").Value = 2119
$ws.Range("C58").Value = "01/04/2022"
$ws.Range("D58").Value = 1956
$ws.Range("C59").Value = "01/07/2022"
$ws.Range("D59").Value = 2041
$ws.Range("C60").Value = "01/10/2022"
$ws.Range("D60").Value = 2100
$ws.Range("C61").Value = "01/01/2023"
$ws.Range("D61").Value = 2111
$ws.Range("C62").Value = "01/04/2023"
$ws.Range("D62").Value = 2159
$ws.Range("C63").Value = "01/07/2023"
$ws.Range("D63").Value = 2063
$ws.Range("C64").Value = "01/10/2023"
$ws.Range("D64").Value = 2050

# Restore the default (un-styled) cell style now that the values are
# safely stored as text, matching the original formatting.
$dateRange.Style = "Normal"
